# Generate Report for Handback
# Updates the localization-status workbook to reflect that the handback
# for "ad2609dc-a353-4226-b0f0-a908752174ce" has completed successfully
# (it was previously flagged as stale / "Ready for handoff").

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the ad2609dc... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn sheet: row 3 is the ad2609dc... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("K3").Value = "2016-09-06 10:37:24"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: row 3 is the ad2609dc... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("K3").Value = "2016-09-06 10:37:42"
$wsDeDe.Range("P3").Value = ""

# The long error-detail text is gone now, so the Error Detail column no
# longer needs to be as wide -- autofit it back down like Excel would
# naturally do when the report is regenerated.
$wsZhCn.Columns.Item(16).AutoFit() | Out-Null
$wsDeDe.Columns.Item(16).AutoFit() | Out-Null
